# Add new test-case columns (FirstName / LastName / ZipCode) next to the
# existing username / password / product table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Entered in the same left-to-right, header/value-interleaved order the
# author used, so the shared-string table gets the same append order:
# FirstName, Syidik, Priam, LastName, ZipCode.
$ws.Range("D1").Value = "FirstName"
$ws.Range("D2").Value = "Syidik"
$ws.Range("E2").Value = "Priam"
$ws.Range("E1").Value = "LastName"
$ws.Range("F1").Value = "ZipCode"
$ws.Range("F2").Value = 13530

# Widen the new columns to match the authored layout.
$ws.Columns.Item(3).ColumnWidth = 25.5
$ws.Columns.Item(4).ColumnWidth = 13.666666666666666
$ws.Columns.Item(5).ColumnWidth = 16

# Move the active selection to F3, as left by the author after data entry.
$ws.Range("F3").Select()
